# Update display name for dtwin id (rename "id" -> "dtwin_id") in the
# semantic_aspect_model_schema and description sheets, widen column A of
# the schema sheet, and refresh the metadata sheet with new commit info.

$wb = $excel.ActiveWorkbook

# --- semantic_aspect_model_schema sheet ---
$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")
$wsSchema.Range("A1").Value = "dtwin_id"
# Target stored column width is 9.6 characters. Excel quantizes ColumnWidth
# to whole pixels (MDW) internally, so the nearest value we can land on is
# requested via 8.8 (which rounds to ~9.67, the closest reachable width).
$wsSchema.Columns.Item(1).ColumnWidth = 8.8

# --- description sheet ---
$wsDescription = $wb.Worksheets.Item("description")
$wsDescription.Range("A5").Value = "dtwin_id"

# --- metadata sheet ---
$wsMetadata = $wb.Worksheets.Item("metadata")
$wsMetadata.Range("B2").Value = "41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B3").Value = "https://github.com/dataspacesolutions/sldt-semantic-models/commit/41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B4").Value = "2025-03-10 14:48:29+00:00"
$wsMetadata.Range("B5").Value = "Adding auto-generated artifacts for new models"
